$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Shift the data block (rows 3:89, cols A:K) up by one row into
#    rows 2:88, carrying both values and styles, so that old row 2
#    (the Hiver/Eté/Année sub-header) is effectively dropped and the
#    data that used to start on row 3 now starts on row 2.
# ---------------------------------------------------------------------
$src = $ws.Range("A3:K89")
$vals = $src.Value()

# Copy formats first (this correctly follows the shift even though
# source and destination overlap).
$src.Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Now copy the values (array assignment also handles the overlap
# correctly, unlike a plain Range.Copy/paste of values).
$ws.Range("A2:K88").Value() = $vals

# The last row no longer has a row below it to pull data from, so it
# becomes fully blank (its previous formatting is dropped too).
$ws.Range("A89:K89").Clear()

# ---------------------------------------------------------------------
# 2. Rewrite row 1 as a single header row.
# ---------------------------------------------------------------------
# Drop whatever formatting the old row 1 / row 2 header cells had
# (e.g. E1 used to carry the "(m3/s)" style) before laying out the
# new headers so no stale style index leaks onto the plain columns.
$ws.Range("A1:K1").Clear()

$ws.Range("A1").Value() = "idx"
$ws.Range("B1").Value() = "idx2"
$ws.Range("C1").Value() = "Name"
$ws.Range("D1").Value() = "Date Start"
$ws.Range("E1").Value() = "Date End"
$ws.Range("F1").Value() = "(m3/s)"
$ws.Range("G1").Value() = "(MW1)"
$ws.Range("H1").Value() = "(MW2)"
$ws.Range("I1").Value() = "(GWh) Winter"
$ws.Range("J1").Value() = "(GWh) Summer"
$ws.Range("K1").Value() = "(GWh) Year"

# F1:K1 use a header style: same font as the rest of the data
# (Arial 9) but without a number format applied to it. Create a
# throwaway named style to mint a fresh cellXf record, apply it, then
# remove the named style again so only the cell format record is left
# behind (matching the workbook's existing convention of not keeping
# custom named cell styles around).
$tmpStyle = $wb.Styles.Add("__tmp_header_style")
$tmpStyle.Font.Name() = "Arial"
$tmpStyle.Font.Size() = 9

$hdr = $ws.Range("F1:K1")
$hdr.Style() = "__tmp_header_style"

$wb.Styles.Item("__tmp_header_style").Delete()

# ---------------------------------------------------------------------
# 3. Match the author's final selection.
# ---------------------------------------------------------------------
$ws.Range("A2:K2").Select()
